# Add the latest S&P 500 earnings-growth data point (2025-03-31, serial 45747,
# value 13.8) to the top of the data table, pushing the existing rows down by
# one row (matching the source sheet, which is sorted most-recent-first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (the first data row),
# shifting all the existing data rows (old row 2 -> new row 3, ..., old row
# 142 -> new row 143) down by one.
$ws.Rows("2:2").Insert()

# The freshly inserted row inherits the header row's formatting by default;
# re-stamp it with the same formatting as the data rows (copy format only
# from the row right below, which still holds the old row 2 formatting).
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Write the new data point into the now-empty row 2.
$ws.Cells.Item(2, 1).Value = 45747
$ws.Cells.Item(2, 2).Value = 13.8
